$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet (keeps column widths/page setup, clears content+format)
$ws.Cells.Clear()

$chf = '"CHF"\ #,##0.00'

# ---- Investitionskosten block ----
$ws.Range("A1").Value = "Investitionskosten"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 255

$ws.Range("A2").Value = "Entwicklungsaufwand = 1 FTE à 14 Monate = CHF 112'800 / 12 * 14"
$ws.Range("B2").Value = 131600
$ws.Range("B2").NumberFormat = $chf

$ws.Range("A3").Value = "Aufbau Entwicklungs- & Testumgebung"
$ws.Range("B3").Value = 35000
$ws.Range("B3").NumberFormat = $chf

$ws.Range("A4").Value = "Initial-Marketing"
$ws.Range("B4").Value = 50000
$ws.Range("B4").NumberFormat = $chf

$ws.Range("A5").Value = "Total Investitionskosten"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Formula = "=SUM(B2:B4)"
$ws.Range("B5").NumberFormat = $chf
$ws.Range("B5").Font.Bold = $true

# ---- Betriebskosten block ----
$ws.Range("A7").Value = "Betriebskosten"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Color = 255
$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").NumberFormat = $chf
$ws.Range("C7").Font.Bold = $true
$ws.Range("C7").NumberFormat = $chf

$ws.Range("A8").Value = "F1 Personalkosten = 1 * 0.5 FTE = CHF 112'800 * 0.5"
$ws.Range("B8").Value = 56400
$ws.Range("B8").NumberFormat = $chf

$ws.Range("A9").Value = "F2 Materialkosten"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = $chf

$ws.Range("A10").Value = "F3 Raumkosten"
$ws.Range("B10").Value = 0
$ws.Range("B10").NumberFormat = $chf

$ws.Range("A11").Value = "F4 Kapitalkosten = 5% Zins auf CHF 700'000"
$ws.Range("B11").Value = 35000
$ws.Range("B11").NumberFormat = $chf

$ws.Range("A12").Value = "F5 Dienstleistungskosten = CHF 800 * 12 Monate"
$ws.Range("B12").Value = 9600
$ws.Range("B12").NumberFormat = $chf

$ws.Range("A13").Value = "F5 Dienstleistungskosten = Werbung"
$ws.Range("B13").Value = 36000
$ws.Range("B13").NumberFormat = $chf
$ws.Range("C13").Font.Bold = $true
$ws.Range("C13").NumberFormat = $chf

$ws.Range("A14").Value = "F6 kalkulatorische Kosten = CHF 700'000 / 3 Jahre = "
$ws.Range("B14").Value = 233333.33
$ws.Range("B14").NumberFormat = $chf

$ws.Range("A15").Value = "F7 Betriebskosten pro Jahr"
$ws.Range("A15").Font.Bold = $true
$ws.Range("B15").Formula = "=SUM(B8:B14)"
$ws.Range("B15").NumberFormat = $chf
$ws.Range("B15").Font.Bold = $true
$ws.Range("C15").Font.Bold = $true
$ws.Range("C15").NumberFormat = $chf

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 59.140625
$ws.Columns.Item(2).ColumnWidth = 21

# ---- Selection / view ----
$ws.Range("A13").Select()

$wb.Windows.Item(1).Left = 32820
$wb.Windows.Item(1).Top = 4530
$wb.Windows.Item(1).Width = 28185
$wb.Windows.Item(1).Height = 15240
